# Update the "Förändrad" (Changed) date column (C) for all existing data
# rows (2-528) from 45181 to 45182.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C528").Value = 45182

# Row 528 gains an explicit custom row height (15) in the new file.
$ws.Rows.Item(528).RowHeight = 15

# Append the new data row (529) with the new avverkningsanmälan entry.
$newRow = 529

$ws.Cells.Item($newRow, 1).Value = "A 42735-2023"

$ws.Cells.Item($newRow, 2).Value = 45181
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45182
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($newRow, 5).Value = "HAGFORS"
$ws.Cells.Item($newRow, 6).Value = "Bergvik skog väst AB"

$ws.Cells.Item($newRow, 7).Value = 2.8

for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}

# Column R (18) stays empty but keeps the wrap-text style used throughout
# the rest of the sheet.
$ws.Cells.Item($newRow, 18).WrapText = $true
